# Edit script implementing the diff:
#   1. Split the "That brings us to today..." paragraph (w14:paraId="3B1C701B")
#      right after "...more colorful." into two paragraphs:
#        - the original two sentences/runs stay in the first paragraph
#        - a new paragraph is added with "Lastly I added an empty game
#          object..." text, followed by the relocated _GoBack bookmark and a
#          trailing run containing a single space
#   2. Move the <w:lastRenderedPageBreak/> marker from the picture run to the
#      "I will quickly go over..." run that now renders as the last line of
#      the previous page.

$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# Step 1: locate + split the target paragraph.
# ---------------------------------------------------------------------
$splitIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $txt = $d.Paragraphs.Item($i).Range.Text
    if ($txt -like "*more colorful.*") {
        $splitIndex = $i
    }
}
if ($splitIndex -eq -1) {
    throw "Could not locate the 'more colorful' paragraph"
}

$p = $d.Paragraphs.Item($splitIndex)
$r = $p.Range

$origParaOpen = "<w:p w14:paraId=`"3B1C701B`" w14:textId=`"73D6A94F`" w:rsidR=`"00051A36`" w:rsidRPr=`"00051A36`" w:rsidRDefault=`"00051A36`" w:rsidP=`"009A02DA`">"
$leadRun = "<w:r><w:t xml:space=`"preserve`">That brings us to today where I </w:t></w:r>"
$origRunOpen = "<w:r w:rsidR=`"00126C0F`">"
$origText = "have added a trigger event on 3 switches. Each one when interacted with will destroy a platform overhead which then drops a box for the player to move along the level. I also added a change color feature when they have been interacted with to make the user experience a bit nicer and make the level more colorful."
$newText = "Lastly I added an empty game object into where the player must get the block (the goal). This will act as a trigger and has a very simple trigger check script attached to it to check to see if the correct block has entered the collider."

$xml = $origParaOpen + $leadRun + $origRunOpen + "<w:t>" + $origText + "</w:t></w:r></w:p><w:p><w:r><w:t>" + $newText + "</w:t></w:r><w:bookmarkStart w:id=`"0`" w:name=`"_GoBack`"/><w:bookmarkEnd w:id=`"0`"/><w:r><w:t xml:space=`"preserve`"> </w:t></w:r></w:p>"
$r.InsertXML($xml)

# ---------------------------------------------------------------------
# Step 2: move <w:lastRenderedPageBreak/> to the correct run.
# ---------------------------------------------------------------------
$targetIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $txt = $d.Paragraphs.Item($i).Range.Text
    if ($txt.StartsWith("I will quickly go over") -and $txt.TrimEnd() -like "*done.") {
        $targetIndex = $i
    }
}
if ($targetIndex -eq -1) {
    throw "Could not locate the 'I will quickly go over...done.' paragraph"
}

$pText = $d.Paragraphs.Item($targetIndex)
$pDrawing = $d.Paragraphs.Item($targetIndex + 1)

$xmlTextPara = "<w:p w14:paraId=`"3EC23BA0`" w14:textId=`"77777777`" w:rsidR=`"006370E8`" w:rsidRDefault=`"006370E8`" w:rsidP=`"006370E8`"><w:r><w:lastRenderedPageBreak/><w:t>I will quickly go over my progress up to this point however to try and cover what’s already been done.</w:t></w:r></w:p>"
$xmlDrawingPara = "<w:p w14:paraId=`"3C7649B5`" w14:textId=`"77777777`" w:rsidR=`"006370E8`" w:rsidRDefault=`"006370E8`" w:rsidP=`"006370E8`"><w:r><w:rPr><w:noProof/><w:lang w:eastAsia=`"en-US`"/></w:rPr><w:drawing><wp:inline distT=`"0`" distB=`"0`" distL=`"0`" distR=`"0`" wp14:anchorId=`"2C29B5F8`" wp14:editId=`"5736B4AC`"><wp:extent cx=`"5486400`" cy=`"3657600`"/><wp:effectExtent l=`"0`" t=`"0`" r=`"0`" b=`"0`"/><wp:docPr id=`"1`" name=`"Picture 1`" descr=`"Red, yellow, and brown boulders on a beach in bright sunshine with a blue sky.`"/><wp:cNvGraphicFramePr><a:graphicFrameLocks xmlns:a=`"http://schemas.openxmlformats.org/drawingml/2006/main`" noChangeAspect=`"1`"/></wp:cNvGraphicFramePr><a:graphic xmlns:a=`"http://schemas.openxmlformats.org/drawingml/2006/main`"><a:graphicData uri=`"http://schemas.openxmlformats.org/drawingml/2006/picture`"><pic:pic xmlns:pic=`"http://schemas.openxmlformats.org/drawingml/2006/picture`"><pic:nvPicPr><pic:cNvPr id=`"2`" name=`"10002048_96.jpg`"/><pic:cNvPicPr/></pic:nvPicPr><pic:blipFill><a:blip r:embed=`"rId7`"><a:extLst><a:ext uri=`"{28A0092B-C50C-407E-A947-70E740481C1C}`"><a14:useLocalDpi xmlns:a14=`"http://schemas.microsoft.com/office/drawing/2010/main`" val=`"0`"/></a:ext></a:extLst></a:blip><a:stretch><a:fillRect/></a:stretch></pic:blipFill><pic:spPr><a:xfrm><a:off x=`"0`" y=`"0`"/><a:ext cx=`"5486400`" cy=`"3657600`"/></a:xfrm><a:prstGeom prst=`"rect`"><a:avLst/></a:prstGeom></pic:spPr></pic:pic></a:graphicData></a:graphic></wp:inline></w:drawing></w:r></w:p>"

$pText.Range.InsertXML($xmlTextPara)
$pDrawing.Range.InsertXML($xmlDrawingPara)

Write-Host "Edit complete."
